$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H12").Value = 140.18182
$ws.Range("I12").Value = 148.88889
$ws.Range("J12").Value = 101
$ws.Range("K12").Value = 148.88889
$ws.Range("L12").Value = 101
$ws.Range("M12").Value = 21.11111
$ws.Range("N12").Value = -441
$ws.Range("H33").Value = 363.2857
$ws.Range("I33").Value = 366.45
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 366.45
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -137.45
$ws.Range("N33").Value = -758
$ws.Range("H112").Value = 1149.7872
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1194.091
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 3582.273
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -5798.272999999999
$ws.Range("H113").Value = 2662.652
$ws.Range("I113").Value = 2480.9333
$ws.Range("J113").Value = 3003.375
$ws.Range("K113").Value = 2480.9333
$ws.Range("L113").Value = 3003.375
$ws.Range("M113").Value = 773.0666999999999
$ws.Range("N113").Value = -9511.375
$ws.Range("H132").Value = 4152.9165
$ws.Range("I132").Value = 3591.0605
$ws.Range("J132").Value = 10333.333
$ws.Range("K132").Value = 10773.1815
$ws.Range("L132").Value = 30999.999
$ws.Range("M132").Value = -8243.181500000001
$ws.Range("N132").Value = -36059.999
$ws.Range("H135").Value = 875.8378
$ws.Range("I135").Value = 754.74286
$ws.Range("J135").Value = 2995
$ws.Range("K135").Value = 6792.68574
$ws.Range("L135").Value = 26955
$ws.Range("M135").Value = -4257.68574
$ws.Range("N135").Value = -32025

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 8951.308999999999
$ws.Range("I32").Value = 9506.387000000001
$ws.Range("J32").Value = 7140
$ws.Range("K32").Value = 9506.387000000001
$ws.Range("L32").Value = 7140
$ws.Range("M32").Value = -9219.387000000001
$ws.Range("N32").Value = -7714
$ws.Range("H61").Value = 11906319
$ws.Range("I61").Value = 14287421
$ws.Range("J61").Value = 809
$ws.Range("K61").Value = 14287421
$ws.Range("L61").Value = 809
$ws.Range("M61").Value = -14287209
$ws.Range("N61").Value = -1233
$ws.Range("H74").Value = 9435486
$ws.Range("I74").Value = 12822044
$ws.Range("J74").Value = 1501.7142
$ws.Range("K74").Value = 12822044
$ws.Range("L74").Value = 1501.7142
$ws.Range("M74").Value = -12821170
$ws.Range("N74").Value = -3249.7142
$ws.Range("H77").Value = 9435486
$ws.Range("I77").Value = 12822044
$ws.Range("J77").Value = 1501.7142
$ws.Range("K77").Value = 64110220
$ws.Range("L77").Value = 7508.571
$ws.Range("M77").Value = -64105852
$ws.Range("N77").Value = -16244.571
$ws.Range("H136").Value = 11906319
$ws.Range("I136").Value = 14287421
$ws.Range("J136").Value = 809
$ws.Range("K136").Value = 42862263
$ws.Range("L136").Value = 2427
$ws.Range("M136").Value = -42859713
$ws.Range("N136").Value = -7527

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 1862.9
$ws.Range("I99").Value = 1847.6666
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1847.6666
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -349.6666
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 2544.7964
$ws.Range("I105").Value = 1486.6786
$ws.Range("J105").Value = 3684.3076
$ws.Range("K105").Value = 1486.6786
$ws.Range("L105").Value = 3684.3076
$ws.Range("M105").Value = 260.3214
$ws.Range("N105").Value = -7178.3076
$ws.Range("H134").Value = 2165.554
$ws.Range("I134").Value = 1392.2245
$ws.Range("K134").Value = 4176.6735
$ws.Range("M134").Value = -1641.6735

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H22").Value = 441.0909
$ws.Range("I22").Value = 241.66667
$ws.Range("J22").Value = 680.4
$ws.Range("K22").Value = 241.66667
$ws.Range("L22").Value = 680.4
$ws.Range("M22").Value = 108.33333
$ws.Range("N22").Value = -1380.4
$ws.Range("H58").Value = 1447.7435
$ws.Range("I58").Value = 842.6286
$ws.Range("K58").Value = 842.6286
$ws.Range("M58").Value = -639.6286
$ws.Range("H99").Value = 1332.6666
$ws.Range("I99").Value = 1178.2106
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1178.2106
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = 319.7893999999999
$ws.Range("N99").Value = -5796
$ws.Range("H105").Value = 1700
$ws.Range("I105").Value = 1700
$ws.Range("K105").Value = 1700
$ws.Range("M105").Value = 47
$ws.Range("H126").Value = 1332.6666
$ws.Range("I126").Value = 1178.2106
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 3534.6318
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -1064.6318
$ws.Range("N126").Value = -13340
$ws.Range("H132").Value = 15626905
$ws.Range("I132").Value = 16668234
$ws.Range("J132").Value = 6957
$ws.Range("K132").Value = 50004702
$ws.Range("L132").Value = 20871
$ws.Range("M132").Value = -50002172
$ws.Range("N132").Value = -25931
$ws.Range("H134").Value = 1557.3959
$ws.Range("I134").Value = 1435.0444
$ws.Range("J134").Value = 3392.6667
$ws.Range("K134").Value = 4305.1332
$ws.Range("L134").Value = 10178.0001
$ws.Range("M134").Value = -1770.1332
$ws.Range("N134").Value = -15248.0001
$ws.Range("H136").Value = 1447.7435
$ws.Range("I136").Value = 842.6286
$ws.Range("K136").Value = 2527.8858
$ws.Range("M136").Value = 22.11419999999998

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H70").Value = 5245
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5245
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15735
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -16365
$ws.Range("H73").Value = 5245
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5245
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15735
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -17919

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H122").Value = 4773.4873
$ws.Range("I122").Value = 4575.615
$ws.Range("J122").Value = 5169.231
$ws.Range("K122").Value = 13726.845
$ws.Range("L122").Value = 15507.693
$ws.Range("M122").Value = -11276.845
$ws.Range("N122").Value = -20407.693
$ws.Range("H136").Value = 10419968
$ws.Range("I136").Value = 12196250
$ws.Range("J136").Value = 16029.286
$ws.Range("K136").Value = 36588750
$ws.Range("L136").Value = 48087.858
$ws.Range("M136").Value = -36586200
$ws.Range("N136").Value = -53187.858

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 1221.3429
$ws.Range("I132").Value = 1037.4364
$ws.Range("J132").Value = 1895.6666
$ws.Range("K132").Value = 3112.3092
$ws.Range("L132").Value = 5686.9998
$ws.Range("M132").Value = -582.3092000000001
$ws.Range("N132").Value = -10746.9998
